$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # Assigning a numeric-looking string straight to .Value lets Excel
    # coerce it to a real number (e.g. "292.50" -> 292.5), which loses
    # the original formatting. Routing it through a text formula and
    # then Paste-Special-Values keeps it a genuine text cell, matching
    # what a user re-typing these values in Excel would end up with.
    $cell = $ws.Range($addr)
    $cell.Formula = '="' + $text.Replace('"', '""') + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

Set-TextValue "D2" '39.801.94'
Set-TextValue "E2" '  +1.44%  '

Set-TextValue "D3" '2.216.86'
Set-TextValue "E3" '  +1.38%  '

Set-TextValue "E4" '  +0.00%  '

Set-TextValue "D5" '292.50'
Set-TextValue "E5" '  -1.17%  '

Set-TextValue "D6" '85.94'
Set-TextValue "E6" '  +5.49%  '

Set-TextValue "D7" '0.514'
Set-TextValue "E7" '  +1.06%  '

Set-TextValue "E8" '  -0.05%  '

Set-TextValue "D9" '0.471'
Set-TextValue "E9" '  +1.25%  '

Set-TextValue "D10" '30.63'
Set-TextValue "E10" '  +5.46%  '

Set-TextValue "D11" '0.0785'
Set-TextValue "E11" '  +2.22%  '

Set-TextValue "D12" '47.43'
Set-TextValue "E12" '  +0.56%  '

Set-TextValue "E13" '  +1.52%  '

Set-TextValue "D14" '6.34'
Set-TextValue "E14" '  +1.74%  '

Set-TextValue "D15" '2.561.21'
Set-TextValue "E15" '  +1.27%  '

Set-TextValue "D16" '14.02'
Set-TextValue "E16" '  +0.81%  '

Set-TextValue "D17" '2.213.61'
Set-TextValue "E17" '  +1.66%  '

Set-TextValue "D18" '0.730'
Set-TextValue "E18" '  +3.31%  '

Set-TextValue "D19" '39.774.73'
Set-TextValue "E19" '  +1.70%  '

Set-TextValue "D20" '0.0₃0881'
Set-TextValue "E20" '  +1.48%  '

Set-TextValue "D21" '11.14'
Set-TextValue "E21" '  +8.58%  '

Set-TextValue "D22" '5.80'

Set-TextValue "D23" '65.43'
Set-TextValue "E23" '  +0.99%  '

Set-TextValue "D24" '235.47'
Set-TextValue "E24" '  +4.62%  '

Set-TextValue "E25" '  -0.07%  '

Set-TextValue "D26" '2.46'
Set-TextValue "E26" '  +2.75%  '

Set-TextValue "D27" '1.83'
Set-TextValue "E27" '  +2.10%  '

Set-TextValue "D28" '22.72'
Set-TextValue "E28" '  +1.06%  '

Set-TextValue "D29" '2.20'
Set-TextValue "E29" '  +1.60%  '

Set-TextValue "D30" '9.22'
Set-TextValue "E30" '  +2.02%  '

Set-TextValue "D31" '32.81'
Set-TextValue "E31" '  +4.34%  '

Set-TextValue "D32" '151.67'
Set-TextValue "E32" '  +1.29%  '

Set-TextValue "E33" '  -0.14%  '

Set-TextValue "D34" '4.93'
Set-TextValue "E34" '  +2.87%  '

Set-TextValue "D35" '0.0718'
Set-TextValue "E35" '  +4.00%  '

Set-TextValue "E36" '  +1.82%  '

Set-TextValue "E37" '  +7.19%  '

Set-TextValue "E38" '  +2.18%  '

Set-TextValue "D39" '15.87'

Set-TextValue "D40" '0.0991'
Set-TextValue "E40" '  +3.24%  '

Set-TextValue "E41" '  +4.40%  '

Set-TextValue "D42" '3.78'
Set-TextValue "E42" '  +5.43%  '

Set-TextValue "D43" '2.062.88'
Set-TextValue "E43" '  +9.28%  '

Set-TextValue "D44" '0.0267'
Set-TextValue "E44" '  +3.49%  '

Set-TextValue "E45" '  +0.42%  '

Set-TextValue "B46" 'EnergySwap'
Set-TextValue "C46" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D46" '17.78'
Set-TextValue "E46" '  +10.63%  '

Set-TextValue "B47" 'FraxShare'
Set-TextValue "C47" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D47" '9.91'
Set-TextValue "E47" '  +11.08%  '

Set-TextValue "D48" '2.59'
Set-TextValue "E48" '  -0.14%  '

Set-TextValue "D49" '2.433.25'
Set-TextValue "E49" '  +1.38%  '

Set-TextValue "D50" '71.26'
Set-TextValue "E50" '  +0.69%  '

Set-TextValue "D51" '88.88'
Set-TextValue "E51" '  +2.68%  '

$excel.CutCopyMode = 0
